$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CasesTab): append ORDER BY clause to the Case ID query
$ws.Range("B2").Value = "MATCH (ss:study_subject)`r`nMATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`r`nWITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files`r`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`r`nMATCH (ss)<-[:sf_of_study_subject]-(sf)`r`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`r`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`r`nMATCH (ss)<-[:demographic_of_study_subject]-(demo)`r`n WHERE ss.disease_subtype IN [`"Paget's Disease`"] `r`nreturn ss.study_subject_id as ``Case ID``,`r`n       p.program_acronym as ``Program Code``,`r`n        p.program_id as Program_ID,`r`n       s.study_acronym as ``Arm``,`r`n       ss.disease_subtype as ``Diagnosis``,`r`n       sf.grouped_recurrence_score AS ``Recurrence Score``,`r`n       d.tumor_size_group AS ``tumor_size``,`r`n       d.er_status AS ``ER Status``,`r`n       d.pr_status AS ``PR Status``,`r`n       demo.age_at_index AS ``Age (years)``,`r`ndemo.survival_time AS ``Survival (days)```r`n order By ss.study_subject_id ASC LIMIT 100"

# Row 3 (SamplesTab): append ORDER BY clause to the Sample ID query
$ws.Range("B3").Value = "MATCH (ss:study_subject)`r`nWITH COLLECT(ss.study_subject_id) AS all_subjects`r`nMATCH (samp:sample)`r`nMATCH (samp)-[:sample_of_study_subject]->(ss)`r`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`r`nMATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`r`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`r`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`r`n WHERE ss.disease_subtype IN [`"Paget's Disease`"] `r`nWITH`r`n    distinct lp,`r`n    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,`r`n    collect(distinct f.file_id) AS files,`r`n    samp, ss, s, p, all_subjects`r`nRETURN`r`n samp.sample_id AS ``Sample ID``,`r`n            ss.study_subject_id AS ``Case ID``,`r`n            p.program_acronym AS ``Program Code``,`r`n            s.study_acronym AS ``Arm``,`r`n            ss.disease_subtype AS ``Diagnosis``,`r`n            samp.tissue_type AS ``Tissue Type``,`r`n            samp.composition AS ``Tissue Composition``,`r`n            samp.sample_anatomic_site AS ``Sample Anatomic Site``,`r`n            samp.method_of_sample_procurement AS ``Sample Procurement Method```r`n order By samp.sample_id ASC LIMIT 100"

# Row 4 (FilesTab): replace the old "order by" clause with the new ORDER BY + LIMIT clause
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`r`nMATCH (f)-[:file_of_sample]->(samp)`r`nMATCH (samp)-[:sample_of_study_subject]->(ss)`r`nMATCH (ss)-[:study_subject_of_study]->(s)`r`nMATCH (s)-[:study_of_program]->(p)`r`nMATCH (d)-[:diagnosis_of_study_subject]->(ss)`r`nMATCH (tp)-[:tp_of_diagnosis]->(d)`r`n WHERE ss.disease_subtype IN [`"Paget's Disease`"] `r`nWITH`r`n        f, parent,p, ss, d,tp, s, samp,`r`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`r`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`r`n        2 as precision`r`nWITH`r`n        f, parent,p, ss, d,tp, s, samp,`r`n        f.file_size /(1024^i) AS value,`r`n        10^precision AS factor,`r`n        units[i] as unit`r`nWITH`r`n        f, parent,p, ss, d,tp, s, samp, unit,`r`n        round(factor * value)/factor AS size`r`nRETURN Distinct`r`n    f.file_name AS ``File Name``,`r`n    head(labels(samp)) AS ``Association``,`r`n    f.file_description AS ``Description``,`r`n    f.file_format AS ``File Format``,`r`n     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`r`n    p.program_acronym AS ``Program Code``,`r`n    s.study_acronym AS ``Arm``,`r`n    ss.study_subject_id AS ``Case ID``,`r`n    samp.sample_id AS ``Sample ID```r`n     order By f.file_name ASC LIMIT 100"

# Row heights grew by one wrapped line for rows 2 and 3 (row 4 unchanged, already at cap)
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# Update selection / scroll position to match the edited file
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
